# Generate Report for Handback
# Marks the zh-cn / de-de handback rows as complete: updates the Status
# column (and the Overview rollup columns that mirror it), records the
# generated target + handback-xliff file names/timestamps, and widens a
# couple of columns that now hold the longer status text / file names.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6e35fa964110a17e7d2166f5c37226707d7e4369/e2e/a.md"

# ---------------------------------------------------------------------
# Overview sheet: status rollup columns (zh-cn / de-de) for both rows,
# plus the two columns are widened to fit the new, longer status text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.1666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.1666666666667

# ---------------------------------------------------------------------
# Per-language detail sheets: zh-cn and de-de share the same shape, only
# the generated xliff file name / handback timestamp differ.
# ---------------------------------------------------------------------
function Set-HandbackRow($ws, $targetXlf, $handbackTime) {
    # Status column
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Target File (I) - now points at the source doc, same as the
    # handoff file link, and is rendered as a hyperlink like column A.
    $ws.Range("I2").Value = "a.md"
    $ws.Range("I3").Value = "a.md"
    $ws.Hyperlinks.Add($ws.Range("I2"), $aMdUrl, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $aMdUrl, "", "", "a.md")
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("I2").Font.Color = 15570276
    $ws.Range("I3").Font.Underline = 2
    $ws.Range("I3").Font.Color = 15570276

    # Latest Handback File (J) - the generated handback xliff.
    $ws.Range("J2").Value = $targetXlf
    $ws.Range("J3").Value = $targetXlf

    # Latest Handback DateTime (K) - when the handback report was made.
    $ws.Range("K2").Value = $handbackTime
    $ws.Range("K3").Value = $handbackTime

    # Column C (Status) / J (Latest Handback File) need to be widened to
    # fit the new values.
    $ws.Columns.Item(3).ColumnWidth = 29.1666666666667
    $ws.Columns.Item(10).ColumnWidth = 39.1666666666667
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HandbackRow $wsZhCn "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-08-17 06:32:34"

$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HandbackRow $wsDeDe "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-08-17 06:32:41"
